$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing designator lists / quantities to account for new parts ---
# (leading apostrophes preserve the existing quote-prefix text style on these cells)

# Row 2: 100nF capacitor group gains C72
$ws.Range("B2").Value = "'C1, C2, C3, C72, C78, C84"
$ws.Range("E2").Value = 6

# Row 16: DMN62D0UW MOSFET gains Q2
$ws.Range("B16").Value = "'Q1, Q2"
$ws.Range("E16").Value = 2

# Row 17: 100 ohm resistor group gains R62
$ws.Range("B17").Value = "'R1, R2, R3, R4, R5, R6, R7, R8, R9, R10, R11, R58, R62"
$ws.Range("E17").Value = 13

# Row 18: 4.7K resistor group gains R59
$ws.Range("B18").Value = "'R12, R13, R14, R15, R16, R59"
$ws.Range("E18").Value = 6

# --- Append new row 37 for the 25MHz oscillator (Y1) ---
# Copy formatting from the last existing data row (36) so borders/styles match.
$ws.Range("A36:F36").Copy()
$ws.Range("A37:F37").PasteSpecial(-4122) # xlPasteFormats

# Leading apostrophes force the text quote-prefix style (matches other text cells in the sheet).
$ws.Range("A37").Value = "'25MHz"
$ws.Range("B37").Value = "'Y1"
$ws.Range("C37").Value = "'ECS_TXO-2016"
$ws.Range("D37").Value = "'OSC 25MHz ECS-TXO-2016"
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = "'XC3163CT-ND"

$excel.CutCopyMode = 0
